$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number (e.g. "213.40") would be
# auto-converted to a numeric cell by a normal Value assignment, same as real
# Excel typing. Force these to stay text (matching the source inlineStr/shared
# string cells) by switching to a Text number format for the assignment, then
# resetting the cell style back to "Normal" so no stray formatting is left behind.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range("D2").Value = '27.930.51'
$ws.Range("E2").Value = '  +1.34%  '
$ws.Range("D3").Value = '1.643.01'
$ws.Range("E3").Value = '  +1.15%  '
$ws.Range("E4").Value = '  +0.03%  '
Set-TextValue $ws.Range("D5") '213.40'
$ws.Range("E5").Value = '  +0.81%  '
$ws.Range("E6").Value = '  -0.27%  '
$ws.Range("E7").Value = '  -0.01%  '
Set-TextValue $ws.Range("D8") '23.73'
$ws.Range("E8").Value = '  +2.29%  '
$ws.Range("E9").Value = '  +0.66%  '
Set-TextValue $ws.Range("D10") '0.0616'
$ws.Range("E10").Value = '  +0.81%  '
$ws.Range("E11").Value = '  -1.82%  '
$ws.Range("D12").Value = '1.876.47'
$ws.Range("E12").Value = '  +1.15%  '
$ws.Range("D13").Value = '1.640.92'
$ws.Range("E13").Value = '  +1.45%  '
$ws.Range("E14").Value = '  +0.81%  '
$ws.Range("E15").Value = '  +3.90%  '
$ws.Range("E16").Value = '  +0.93%  '
$ws.Range("D17").Value = '27.912.31'
$ws.Range("E17").Value = '  +1.39%  '
Set-TextValue $ws.Range("D18") '230.14'
$ws.Range("E18").Value = '  -0.63%  '
$ws.Range("E19").Value = '  +0.67%  '
Set-TextValue $ws.Range("D20") '7.64'
$ws.Range("E20").Value = '  +1.33%  '
$ws.Range("E21").Value = '  -0.02%  '
Set-TextValue $ws.Range("D22") '10.98'
$ws.Range("E22").Value = '  +5.15%  '
$ws.Range("E23").Value = '  +1.49%  '
$ws.Range("E24").Value = '  +2.70%  '
Set-TextValue $ws.Range("D25") '152.15'
$ws.Range("E25").Value = '  +1.97%  '
Set-TextValue $ws.Range("D26") '6.92'
$ws.Range("E26").Value = '  +0.52%  '
$ws.Range("E27").Value = '  +0.85%  '
$ws.Range("E28").Value = '  +1.06%  '
$ws.Range("E29").Value = '  +0.03%  '
$ws.Range("E30").Value = '  +1.11%  '
$ws.Range("E31").Value = '  +0.23%  '
Set-TextValue $ws.Range("D32") '3.33'
$ws.Range("E32").Value = '  +1.79%  '
$ws.Range("D33").Value = '1.424.38'
$ws.Range("E33").Value = '  -3.06%  '
$ws.Range("E34").Value = '  +1.08%  '
Set-TextValue $ws.Range("D35") '1.57'
$ws.Range("E35").Value = '  +1.79%  '
$ws.Range("E36").Value = '  -0.03%  '
Set-TextValue $ws.Range("D37") '0.889'
$ws.Range("E37").Value = '  +1.77%  '
$ws.Range("E38").Value = '  +1.10%  '
$ws.Range("E39").Value = '  -1.02%  '
Set-TextValue $ws.Range("D40") '0.558'
$ws.Range("E40").Value = '  +0.21%  '
$ws.Range("E41").Value = '  +2.64%  '
Set-TextValue $ws.Range("D43") '67.66'
$ws.Range("E43").Value = '  +0.05%  '
$ws.Range("E44").Value = '  +1.06%  '
$ws.Range("E45").Value = '  +2.71%  '
$ws.Range("E46").Value = '  +2.93%  '
$ws.Range("E47").Value = '  -0.06%  '
$ws.Range("D48").Value = '1.784.92'
$ws.Range("E48").Value = '  +1.19%  '
Set-TextValue $ws.Range("D49") '88.77'
$ws.Range("E49").Value = '  +1.44%  '
$ws.Range("E50").Value = '  +0.72%  '
$ws.Range("E51").Value = '  +0.29%  '
